# docs/resources/Mappings.xlsx -- "chore: docs and templates update (#4243)"
# update to fit with catalogue data model >4.x
#
# Sheet1 = "Dataset mappings" (A1:G1) keeps the same header text, only its
# selection + column widths change.
# Sheet2 = "Variable mappings" (A1:M1) renames column I "match" -> "repeats"
# and column J "status" -> "match"; selection + column widths also change.

$wb = $excel.ActiveWorkbook

$dsSheet  = $wb.Worksheets.Item("Dataset mappings")
$varSheet = $wb.Worksheets.Item("Variable mappings")

# ----- Variable mappings: update the two header cells that actually changed
$varSheet.Range("I1").Value = "repeats"
$varSheet.Range("J1").Value = "match"

# ----- Column widths (best-fit-style widths on both sheets)
$dsSheet.Range("A1").ColumnWidth = 5.166666666666667
$dsSheet.Range("B1").ColumnWidth = 11.5
$dsSheet.Range("C1").ColumnWidth = 4.5
$dsSheet.Range("D1").ColumnWidth = 10.666666666666666
$dsSheet.Range("E1").ColumnWidth = 4.0
$dsSheet.Range("F1").ColumnWidth = 8.833333333333334
$dsSheet.Range("G1").ColumnWidth = 4.833333333333333

$varSheet.Range("A1").ColumnWidth = 5.166666666666667
$varSheet.Range("B1").ColumnWidth = 11.5
$varSheet.Range("C1").ColumnWidth = 12.666666666666666
$varSheet.Range("D1").ColumnWidth = 30.833333333333332
$varSheet.Range("E1").ColumnWidth = 29.333333333333332
$varSheet.Range("F1").ColumnWidth = 4.5
$varSheet.Range("G1").ColumnWidth = 10.666666666666666
$varSheet.Range("H1").ColumnWidth = 10.833333333333334
$varSheet.Range("I1").ColumnWidth = 5.833333333333333
$varSheet.Range("J1").ColumnWidth = 4.833333333333333
$varSheet.Range("K1").ColumnWidth = 8.833333333333334
$varSheet.Range("L1").ColumnWidth = 4.833333333333333
$varSheet.Range("M1").ColumnWidth = 8.333333333333334

# ----- Selection (cosmetic, matches the saved cursor position in the file)
$dsSheet.Range("B4").Select() | Out-Null
$varSheet.Range("C4").Select() | Out-Null

# Re-select the first sheet so it stays the active/tab-selected sheet.
$dsSheet.Select() | Out-Null
